$wb = $excel.ActiveWorkbook

# --- tc0002: add a new step row (duplicate-id repro: q7 / asdf / qewr / zxcv) ---
$wsTc0002 = $wb.Worksheets.Item("tc0002")
$wsTc0002.Cells.Item(8, 1).Value = "q7"
$wsTc0002.Cells.Item(8, 2).Value = "asdf"
$wsTc0002.Cells.Item(8, 3).Value = "qewr"
$wsTc0002.Cells.Item(8, 4).Value = "zxcv"

# --- asdf: add a new step row (q6 / asdf / asdf / asdf), formatted like row 7 ---
$wsAsdf = $wb.Worksheets.Item("asdf")
[void]$wsAsdf.Rows("7").Copy()
[void]$wsAsdf.Rows("8").Insert()
$wsAsdf.Cells.Item(8, 1).Value = "q6"
$wsAsdf.Cells.Item(8, 2).Value = "asdf"
$wsAsdf.Cells.Item(8, 3).Value = "asdf"
$wsAsdf.Cells.Item(8, 4).Value = "asdf"
[void]$wsAsdf.Range("D8").Select()

# --- 827asaf: add a new step row (q6 / asdf / asdf / qwer), formatted like row 7 ---
$ws827 = $wb.Worksheets.Item("827asaf")
[void]$ws827.Rows("7").Copy()
[void]$ws827.Rows("8").Insert()
$ws827.Cells.Item(8, 1).Value = "q6"
$ws827.Cells.Item(8, 2).Value = "asdf"
$ws827.Cells.Item(8, 3).Value = "asdf"
$ws827.Cells.Item(8, 4).Value = "qwer"
[void]$ws827.Range("D8").Select()

# --- tc0002 ends up the active sheet/tab, with D8 selected ---
$wsTc0002.Activate()
[void]$wsTc0002.Range("D8").Select()
